$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values for specific rows, per repull of data.
$updates = @{
    8  = 5
    9  = 0
    22 = 1
    25 = 2
    27 = 1
    29 = 2
    33 = 0
    34 = -2
    41 = -2
    48 = 3
    50 = 3
    51 = -4
    52 = -2
    53 = 2
    55 = -4
    58 = -1
    67 = 3
    71 = -3
    74 = -1
    77 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
